# Applies the "frames capability (WIP)" edit:
#   - remove the blank paragraph between the Speech-bubble hyperlink and
#     the "Thought bubble sticker:" line
#   - turn the (now) blank paragraph after "Thought bubble sticker:" into
#     a hyperlink to the pngall thought-bubble PNG
#   - add a hyperlink to the onlygfx explosion-bubble PNG (plus a trailing
#     space) into the last paragraph, in front of the _GoBack bookmark

$d = $word.ActiveDocument

# --- 1. Drop the stray empty paragraph right before "Thought bubble sticker:" ---
$blank = $d.Paragraphs.Item(4)
$blank.Range.Delete()

# After the delete the document collapses to 7 paragraphs:
#   1 List of graphical assets found from the internet:
#   2 Speech bubble sticker:
#   3 <hyperlink> http://pluspng.com/png-143720.html
#   4 Thought bubble sticker:
#   5 (blank - to become the pngall hyperlink)
#   6 Exclamation sticker:
#   7 (bookmark paragraph - gets the onlygfx hyperlink + space)

# --- 2. Turn the blank paragraph after "Thought bubble sticker:" into a hyperlink ---
$thoughtUrl = "http://www.pngall.com/thought-bubble-png/download/15113"
$thoughtPara = $d.Paragraphs.Item(5)
$thoughtStart = $thoughtPara.Range.Start

$insertPoint = $d.Range($thoughtStart, $thoughtStart)
$insertPoint.InsertBefore($thoughtUrl)

$thoughtRange = $d.Range($thoughtStart, $thoughtStart + $thoughtUrl.Length)
$d.Hyperlinks.Add($thoughtRange, $thoughtUrl, $null, $null, $thoughtUrl) | Out-Null

# --- 3. Add the onlygfx hyperlink + trailing space before the _GoBack bookmark ---
$explosionUrl = "https://www.onlygfx.com/5-comic-explosion-bubble-png-transparent-svg-vector/"
$bookmarkPara = $d.Paragraphs.Item(7)
$bookmarkStart = $bookmarkPara.Range.Start

$insertPoint2 = $d.Range($bookmarkStart, $bookmarkStart)
$insertPoint2.InsertBefore($explosionUrl + " ")

$explosionRange = $d.Range($bookmarkStart, $bookmarkStart + $explosionUrl.Length)
$d.Hyperlinks.Add($explosionRange, $explosionUrl, $null, $null, $explosionUrl) | Out-Null
